$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Edad" column before the "Salto" column (currently column C),
# shifting the existing Salto/Fuerza/Velocidad/Puntuacion/Tiempo/Total
# columns one place to the right (C:H -> D:I), without doing a real
# column-insert (which would also relocate the custom column-width
# definition on column F).
$lastRow = 11
$lastCol = 8

for ($row = 1; $row -le $lastRow; $row++) {
    for ($col = $lastCol; $col -ge 3; $col--) {
        $val = $ws.Cells.Item($row, $col).Value()
        $ws.Cells.Item($row, $col + 1).Value = $val
    }
}

# Header for the new column
$ws.Cells.Item(1, 3).Value = "Edad"

# Age values for the new column
$ages = @(34, 53, 23, 29, 33, 56, 42, 38, 32, 47)
for ($i = 0; $i -lt $ages.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $ages[$i]
}

# Match the final selection left behind in the saved file
$null = $ws.Range("C11").Select()
